$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.987.16'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '1.827.43'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = '''311.87'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").Value = '''1.004'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D8").Value = '''0.3700'
$ws.Range("E8").Value = '  +1.65%  '
$ws.Range("D9").Value = '''0.07336'
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").Value = '''0.8742'
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("D11").Value = '''0.07938'
$ws.Range("E11").Value = '  +4.20%  '
$ws.Range("D12").Value = '''19.81'
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("D13").Value = '1.835.10'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").Value = '''5.340'
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '''6.546'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '''91.37'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '''0.000008885'
$ws.Range("E18").Value = '  +2.89%  '
$ws.Range("D19").Value = '''1.005'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '''14.77'
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("D21").Value = '27.172.31'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").Value = '''5.113'
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D23").Value = '''10.56'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '2.089.81'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '''153.10'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '''1.849'
$ws.Range("E26").Value = '  -1.41%  '
$ws.Range("D27").Value = '''18.39'
$ws.Range("E27").Value = '  +0.85%  '
$ws.Range("D28").Value = '''2.039'
$ws.Range("E28").Value = '  -2.59%  '
$ws.Range("D29").Value = '''5.145'
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("D30").Value = '''115.44'
$ws.Range("E30").Value = '  -0.66%  '
$ws.Range("D31").Value = '''0.08903'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").Value = '''2.961'
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").Value = '''0.7320'
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").Value = '''4.432'
$ws.Range("D35").Value = '''1.127'
$ws.Range("D36").Value = '''2.474'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '''0.01954'
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").Value = '''0.05226'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("D42").Value = '''0.5166'
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '''8.189'
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("D46").Value = '''1.004'
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '''10.16'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '''102.57'
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("D50").Value = '''0.06197'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '''64.73'
$ws.Range("E51").Value = '  +0.28%  '
